$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 30, pushing existing rows 30..131 down to 32..133
$ws.Rows("30:31").Insert()

# Fill in the new row 30 with new data
$ws.Range("A30").Value = 11
$ws.Range("B30").Value = "Vega Monumental Concepción"
$ws.Range("C30").Value = "Bíobío"
$ws.Range("D30").Value = 45014
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103002
$ws.Range("J30").Value = "Ciruela"
$ws.Range("K30").Value = "Angeleno"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 9500
$ws.Range("Q30").Value = "$/bandeja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 528
$ws.Range("T30").Value = 18

# Fill in the new row 31 with new data
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 45014
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103002
$ws.Range("J31").Value = "Ciruela"
$ws.Range("K31").Value = "Angeleno"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 7000
$ws.Range("O31").Value = 7000
$ws.Range("P31").Value = 7000
$ws.Range("Q31").Value = "$/bandeja 18 kilos granel"
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 389
$ws.Range("T31").Value = 18
